$wb = $excel.ActiveWorkbook

# Regenerate the status report: the handoff step has completed and the
# item has moved from "Ready for handoff" into active translation.
# Update the Status value everywhere it appears: the per-locale columns
# on the Overview sheet, and the Status column on each locale's detail
# sheet.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# The Status column is narrower for the new (shorter) status text -
# shrink the matching columns on all three sheets to keep them in sync.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
